$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in previously-missing value in column D (row 3)
$ws.Cells.Item(3, 4).Value = -14.2

# Clear out value in column D (row 5) - now treated as missing
$ws.Cells.Item(5, 4).ClearContents()

# Fill in previously-missing value in column D (row 21)
$ws.Cells.Item(21, 4).Value = -14.3

# Clear out value in column D (row 23) - now treated as missing
$ws.Cells.Item(23, 4).ClearContents()

# Remove the "RM 232" row entirely (row 26)
$ws.Rows(26).Delete()

# Remove the "SC 92" row entirely (now shifted up to row 27)
$ws.Rows(27).Delete()

# Fill in previously-missing value in column D for "SC 193" row (now row 32 after the two deletions above)
$ws.Cells.Item(32, 4).Value = -14.7
